$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Force the value to be stored as text (matching the original inline-string
    # cell type) rather than letting Excel auto-convert numeric-looking strings
    # (e.g. "6.19", "0.999") into floating point numbers. The leading
    # apostrophe forces text entry; resetting the style back to "Normal"
    # afterwards removes the transient quote-prefix formatting so the cell's
    # style stays identical to its original (unstyled) state.
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "58.589.50"
Set-TextValue $ws.Range("E2") "  -3.95%  "
Set-TextValue $ws.Range("D3") "2.554.20"
Set-TextValue $ws.Range("E3") "  -1.97%  "
Set-TextValue $ws.Range("E4") "  +0.05%  "
Set-TextValue $ws.Range("D5") "505.74"
Set-TextValue $ws.Range("E5") "  -3.45%  "
Set-TextValue $ws.Range("D6") "143.47"
Set-TextValue $ws.Range("E6") "  -7.41%  "
Set-TextValue $ws.Range("D8") "0.554"
Set-TextValue $ws.Range("E8") "  -6.67%  "
Set-TextValue $ws.Range("D9") "2.562.86"
Set-TextValue $ws.Range("E9") "  -2.03%  "
Set-TextValue $ws.Range("D10") "6.19"
Set-TextValue $ws.Range("E10") "  -7.41%  "
Set-TextValue $ws.Range("D11") "0.101"
Set-TextValue $ws.Range("E11") "  -3.85%  "
Set-TextValue $ws.Range("D12") "0.330"
Set-TextValue $ws.Range("E12") "  -4.89%  "
Set-TextValue $ws.Range("E13") "  -1.06%  "
Set-TextValue $ws.Range("D14") "3.003.10"
Set-TextValue $ws.Range("E14") "  -1.92%  "
Set-TextValue $ws.Range("D15") "58.625.13"
Set-TextValue $ws.Range("E15") "  -3.92%  "
Set-TextValue $ws.Range("D16") "20.55"
Set-TextValue $ws.Range("E16") "  -5.53%  "
Set-TextValue $ws.Range("D17") "0.0000134"
Set-TextValue $ws.Range("E17") "  -5.34%  "
Set-TextValue $ws.Range("D18") "2.556.65"
Set-TextValue $ws.Range("E18") "  -2.06%  "
Set-TextValue $ws.Range("D19") "4.51"
Set-TextValue $ws.Range("E19") "  -5.46%  "
Set-TextValue $ws.Range("D20") "332.78"
Set-TextValue $ws.Range("E20") "  -6.12%  "
Set-TextValue $ws.Range("D21") "10.05"
Set-TextValue $ws.Range("E21") "  -5.03%  "
Set-TextValue $ws.Range("D22") "0.996"
Set-TextValue $ws.Range("E22") "  -0.35%  "
Set-TextValue $ws.Range("D23") "5.94"
Set-TextValue $ws.Range("E23") "  -4.61%  "
Set-TextValue $ws.Range("D24") "59.81"
Set-TextValue $ws.Range("E24") "  -1.98%  "
Set-TextValue $ws.Range("D25") "0.406"
Set-TextValue $ws.Range("E25") "  -4.94%  "
Set-TextValue $ws.Range("D26") "0.999"
Set-TextValue $ws.Range("E26") "  +0.17%  "
Set-TextValue $ws.Range("E27") "  -6.00%  "
Set-TextValue $ws.Range("D28") "0.0₃0776"
Set-TextValue $ws.Range("E28") "  -8.60%  "
Set-TextValue $ws.Range("D29") "6.87"
Set-TextValue $ws.Range("E29") "  -7.61%  "
Set-TextValue $ws.Range("E30") "  +0.00%  "
Set-TextValue $ws.Range("B31") "Monero"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D31") "149.14"
Set-TextValue $ws.Range("E31") "  +0.62%  "
Set-TextValue $ws.Range("B32") "Aptos"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D32") "5.84"
Set-TextValue $ws.Range("E32") "  -7.25%  "
Set-TextValue $ws.Range("B33") "EthereumClassic"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D33") "18.52"
Set-TextValue $ws.Range("E33") "  -4.63%  "
Set-TextValue $ws.Range("E34") "  -4.01%  "
Set-TextValue $ws.Range("D35") "3.92"
Set-TextValue $ws.Range("E35") "  -6.67%  "
Set-TextValue $ws.Range("D36") "0.899"
Set-TextValue $ws.Range("E36") "  -3.55%  "
Set-TextValue $ws.Range("E37") "  -8.35%  "
Set-TextValue $ws.Range("D38") "35.96"
Set-TextValue $ws.Range("E38") "  -1.44%  "
Set-TextValue $ws.Range("D39") "0.820"
Set-TextValue $ws.Range("E39") "  -6.68%  "
Set-TextValue $ws.Range("D40") "287.31"
Set-TextValue $ws.Range("E40") "  -1.48%  "
Set-TextValue $ws.Range("E41") "  -8.71%  "
Set-TextValue $ws.Range("D42") "3.52"
Set-TextValue $ws.Range("E42") "  -7.79%  "
Set-TextValue $ws.Range("D43") "0.999"
Set-TextValue $ws.Range("E43") "  +0.20%  "
Set-TextValue $ws.Range("B44") "Mantle"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D44") "0.606"
Set-TextValue $ws.Range("E44") "  -2.84%  "
Set-TextValue $ws.Range("B45") "Stellar"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D45") "0.0980"
Set-TextValue $ws.Range("E45") "  -3.69%  "
Set-TextValue $ws.Range("D46") "0.0532"
Set-TextValue $ws.Range("E46") "  -5.19%  "
Set-TextValue $ws.Range("E47") "  -4.78%  "
Set-TextValue $ws.Range("E48") "  -0.06%  "
Set-TextValue $ws.Range("E49") "  -4.79%  "
Set-TextValue $ws.Range("D50") "4.52"
Set-TextValue $ws.Range("E50") "  -10.37%  "
Set-TextValue $ws.Range("D51") "1.910.62"
Set-TextValue $ws.Range("E51") "  -3.00%  "
